$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.731.29'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.33%  '
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.414.74'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.58%  '
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.22'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.83%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '127.72'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -2.22%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +6.42%  '
# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.09%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.727'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +7.13%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.137'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +5.90%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.33'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.88%  '
# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.141'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.10%  '
# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.06'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +8.80%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.954.75'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.65%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.22'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.63%  '
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +40.03%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.405.67'
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.18'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +5.23%  '
# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.96%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '61.700.37'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.44%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '442.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +42.53%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.15'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +9.82%  '
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.82%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.76%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.22'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.45%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.74'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +11.47%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.62'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.45%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.73'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.68%  '
# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -6.16%  '
# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.32%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.91'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.68%  '
# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.71%  '
# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.11%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.23'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.04%  '
# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.09%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0494'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.10%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.09'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.69%  '
# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.07%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.35'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.38%  '
# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.57%  '
# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.62%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.315'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.87%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.73'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.26%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.17'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +6.56%  '
# Row 45
$ws.Range("B45").Value = 'WEMIXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.53'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +14.51%  '
# Row 46
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.97'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.31%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.47'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.21%  '
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.65%  '
# Row 49
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.766.38'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.90%  '
# Row 50
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.11'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.36%  '
# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.136'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +17.93%  '
